$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update R23 and R27 values
$ws.Range("R23").Value = 2
$ws.Range("R27").Value = 2

# Fill in row 32 (subject 30) - B32 already "F" condition (string index 18), keep
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 4
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = 4
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 4
$ws.Range("I32").Value = 4
$ws.Range("J32").Value = 5
$ws.Range("K32").Value = 5
$ws.Range("L32").Value = 3
$ws.Range("M32").Value = 5
$ws.Range("N32").Value = 4
$ws.Range("O32").Value = 5
$ws.Range("P32").Value = 4
$ws.Range("Q32").Value = 5

# Fill in row 33 (subject 31) - B33 changes from "G" to "F"
$ws.Range("B33").Value = "F"
$ws.Range("C33").Value = 5
$ws.Range("D33").Value = 4
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 2
$ws.Range("I33").Value = 4
$ws.Range("J33").Value = 4
$ws.Range("K33").Value = 3
$ws.Range("L33").Value = 1
$ws.Range("M33").Value = 3
$ws.Range("N33").Value = 3
$ws.Range("O33").Value = 4
$ws.Range("P33").Value = 2
$ws.Range("Q33").Value = 3

# Update sheet view - scroll and selection
$ws.Activate()
$ws.Range("Q34").Select()
$excel.ActiveWindow.ScrollRow = 10
